$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The draw-results sheet is a simple append-only log; today's row goes
# right after the last populated row (row 62 -> new row 63).
$newRow = 63

# Columns A (Date) and C (Phase) hold numeric-looking strings
# ("2025-11-18", "251118"). Force them to be stored as text (matching
# every other row in the sheet, where t="str") instead of letting Excel
# auto-infer a date/number, by pre-setting the cell number format to
# Text ("@") before assigning the value.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 3).NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-11-18"
$ws.Cells.Item($newRow, 2).Value = "Pick 3"
$ws.Cells.Item($newRow, 3).Value = "251118"
$ws.Cells.Item($newRow, 4).Value = "9-6-7"
$ws.Cells.Item($newRow, 5).Value = "2025-11-18T21:40:41.381+04:00"
